# Overview - Fashion Retailers.xlsx
# Commit: Update currency pairs, update $UA price
#
# The "Main" sheet keeps a small set of manual currency-conversion inputs
# (USDGBP, EURGBP, JPYGBP) in E50:E52. Every other changed number in the
# workbook (F3, G3, H3, ... AE41, plus the inverse-rate helpers in F50:F52)
# is a formula that derives from these three cells (directly, or through the
# external-workbook-linked "$UA" style rows), so updating just the inputs
# and letting Excel recalculate reproduces the whole set of value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Make "Main" the active sheet/tab, matching the saved file's UI state.
$ws.Activate()

# USDGBP rate: 0.82 -> 0.80
$ws.Range("E50").Value = 0.8

# EURGBP rate: 0.87 -> 0.88
$ws.Range("E51").Value = 0.88

# JPYGBP rate: 0.006 -> 0.0062
$ws.Range("E52").Value = 0.0062

# Leave the cursor on the cell that was last edited.
$ws.Range("E52").Select()
